# Chore Tracker update: replace the sample placeholder rows with real
# chore/assignee data, clear the date formatting on column C (now holding
# literal "#########" text instead of a date), set every cell's Status
# (column D) to the values from the new data set, and move the active
# selection to F5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New task / assignee / "due date" / status data (rows 2-7)
$data = @(
    @("Laundry",       "Gapper", "#########", $false),
    @("Siege",         "Warb",   "#########", $true),
    @("Dishes",        "Gapper", "#########", $true),
    @("Laundry",       "Gapper", "#########", $true),
    @("Floor",         "Gapper", "#########", $false),
    @("Sample Task 6",  "Warb",  "#########", $false)
)

# Column C used to carry a custom date number format (style index 2). The
# new content is plain text, so clear any existing formatting on that
# column before writing the new values.
$ws.Range("C2:C7").ClearFormats()

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Move the active selection to F5 (was A2:A7)
[void]$ws.Range("F5").Select()
